$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.3325634519964858, -1.577164117715349, 0.1583661481491326)
    3  = @(0.5362673463982073, 2.692348254556006, 0.1009111604271066)
    4  = @(0.4863417279209184, 4.190925488841986, 0.1137652762049108)
    5  = @(0.6631075989225439, 5.443480049424021, 0.1931707690100237)
    6  = @(0.6775057263139524, 6.112545600687255, 0.2101936556005051)
    7  = @(0.443581991011115, 6.436510515104551, 0.08830878254179103)
    8  = @(0.6745182959672391, 6.861214481091074, 0.1942603081987745)
    9  = @(0.5690983528995748, 7.508962356659595, 0.1950218386347838)
    10 = @(0.5805728319597545, 7.541483340287529, 0.1864556023826213)
    11 = @(0.5013123676011131, 8.078763114201099, 0.1046331517116006)
    12 = @(0.5538508295498324, 8.136658617034685, 0.09199700183889786)
    13 = @(0.5443277979225418, 8.4438213234104, 0.1653917755315666)
    14 = @(0.4734461023878559, 8.517659546701067, 0.2215697407280473)
    15 = @(0.7319111827345764, 8.886243262158843, 0.2105858936212124)
    16 = @(0.2719028202929124, 9.086684146926034, 0.1643884560026617)
    17 = @(0.6553337074312454, 9.144604967488755, 0.1538476224249944)
    18 = @(0.3744091181736931, 9.390416267008765, 0.1865534655677196)
    19 = @(0.6905295855619613, 9.472828200603963, 0.1665361906093438)
    20 = @(0.2777390065955421, 9.707266478039442, 0.1690429448834274)
    21 = @(0.5784279659116576, 9.565568645208231, 0.2050538393862699)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}
